$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from G1 (existing "sum" header) onto the new H1 header cell
# so the new column reuses the same cell style as the other headers.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 1
